# Plano de Ação workbook update
# - Rename table "Dados" headers for forecast dates
# - Update task rows (status, week, assignees, task names)
# - Add two new rows at the end of the table
# - Fix previously-broken date formulas for several rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados do plano de Ação")
$lo = $ws.ListObjects.Item("Dados")

# --- Rename the forecast date columns (G/H headers) ---
$ws.Range("G6").Value2 = "Previsão de início"
$ws.Range("H6").Value2 = "Previsão de término"

# --- Add two new rows to the table so it spans B6:H21 ---
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# --- Row 7: Análise de Entregaveis Sprint 1 -> Atualização GitHub do Projeto ---
$ws.Range("B7").Value2 = "Atualização GitHub do Projeto"
$ws.Range("C7").Value2 = "Concluído"
$ws.Range("D7").Value2 = "Semana 2-A"
$ws.Range("E7").Value2 = "Yuri"
$ws.Range("F7").Value2 = "Kaiqui"

# --- Row 8: Documentação Segunda Versão -> Documentação Atualizada ---
$ws.Range("B8").Value2 = "Documentação Atualizada"
$ws.Range("C8").Value2 = "Concluído"
$ws.Range("D8").Value2 = "Semana 2-A"
$ws.Range("E8").Value2 = "Yuri"
$ws.Range("F8").Value2 = "Thalita"

# --- Row 9: Modelagem Banco de Dados Primeira Versão (kept) ---
$ws.Range("B9").Value2 = "Modelagem Banco de Dados Primeira Versão"
$ws.Range("C9").Value2 = "Concluído"
$ws.Range("D9").Value2 = "Semana 2-A"
$ws.Range("E9").Value2 = "Yuri"
$ws.Range("F9").Value2 = "Ester"

# --- Row 10: Calculadora Financeira Segunda Versão, now Em andamento ---
$ws.Range("B10").Value2 = "Calculadora Financeira Segunda Versão"
$ws.Range("C10").Value2 = "Em andamento"
$ws.Range("D10").Value2 = "Semana 2-B"
$ws.Range("E10").Value2 = "Ester"
$ws.Range("F10").Value2 = "Yuri"

# --- Row 11: Prototipo do Site Ajustado, now Em andamento ---
$ws.Range("B11").Value2 = "Prototipo do Site Ajustado"
$ws.Range("C11").Value2 = "Em andamento"
$ws.Range("D11").Value2 = "Semana 2-B"
$ws.Range("E11").Value2 = "Ester"
$ws.Range("F11").Value2 = "Ester, Kaiqui e Paulo"

# --- Row 12: Página Inicial Estática, week moves to Semana 2-C, dates fixed ---
$ws.Range("B12").Value2 = "Página Inicial Estática"
$ws.Range("C12").Value2 = "Não Iniciado"
$ws.Range("D12").Value2 = "Semana 2-C"
$ws.Range("G12").Formula = "=DATE(YEAR(TODAY()),9,9)"
$ws.Range("H12").Formula = "=DATE(YEAR(TODAY()),9,13)"

# --- Row 13: Página Login Estática ---
$ws.Range("B13").Value2 = "Página Login Estática"
$ws.Range("C13").Value2 = "Não Iniciado"
$ws.Range("D13").Value2 = "Semana 2-C"
$ws.Range("G13").Formula = "=DATE(YEAR(TODAY()),9,9)"
$ws.Range("H13").Formula = "=DATE(YEAR(TODAY()),9,13)"

# --- Row 14: Página Cadastro Estática ---
$ws.Range("B14").Value2 = "Página Cadastro Estática"
$ws.Range("C14").Value2 = "Não Iniciado"
$ws.Range("D14").Value2 = "Semana 2-C"
$ws.Range("G14").Formula = "=DATE(YEAR(TODAY()),9,9)"
$ws.Range("H14").Formula = "=DATE(YEAR(TODAY()),9,13)"

# --- Row 15: Página Dashboard Estática ---
$ws.Range("B15").Value2 = "Página Dashboard Estática"
$ws.Range("C15").Value2 = "Não Iniciado"
$ws.Range("D15").Value2 = "Semana 2-C"
$ws.Range("G15").Formula = "=DATE(YEAR(TODAY()),9,9)"
$ws.Range("H15").Formula = "=DATE(YEAR(TODAY()),9,13)"

# --- Row 16: Integração Captura de Dados com ChartJs -> Teste com Sensor do Projeto + Gráficos ---
$ws.Range("B16").Value2 = "Teste com Sensor do Projeto + Gráficos"
$ws.Range("C16").Value2 = "Concluído"
$ws.Range("D16").Value2 = "Semana 2-A"
$ws.Range("E16").Value2 = "Ester"
$ws.Range("F16").Value2 = "Thalita e Yuri"
$ws.Range("G16").Formula = "=DATE(YEAR(TODAY()),9,2)"
$ws.Range("H16").Formula = "=DATE(YEAR(TODAY()),9,6)"

# --- Row 17: Especificação de Métricas/Analytics, week moves to Semana 2-D ---
$ws.Range("B17").Value2 = "Especificação de Métricas/Analytics"
$ws.Range("C17").Value2 = "Não Iniciado"
$ws.Range("D17").Value2 = "Semana 2-D"
$ws.Range("G17").Formula = "=DATE(YEAR(TODAY()),9,16)"
$ws.Range("H17").Formula = "=DATE(YEAR(TODAY()),9,20)"

# --- Row 18: Diagrama de solução ---
$ws.Range("B18").Value2 = "Diagrama de solução"
$ws.Range("C18").Value2 = "Não Iniciado"
$ws.Range("D18").Value2 = "Semana 2-D"
$ws.Range("G18").Formula = "=DATE(YEAR(TODAY()),9,16)"
$ws.Range("H18").Formula = "=DATE(YEAR(TODAY()),9,20)"

# --- Row 19: Script Banco de Dados -> Atividades organizadas na ferramenta de Gestão ---
$ws.Range("B19").Value2 = "Atividades organizadas na ferramenta de Gestão (Sprints / Atividades)"
$ws.Range("C19").Value2 = "Concluído"
$ws.Range("D19").Value2 = "Semana 2-B"
$ws.Range("E19").Value2 = "Ester"
$ws.Range("F19").Value2 = "Ester"

# --- Row 20 (new): Planilha de Riscos do Projeto ---
$ws.Range("B20").Value2 = "Planilha de Riscos do Projeto"
$ws.Range("C20").Value2 = "Concluído"
$ws.Range("D20").Value2 = "Semana 2-B"
$ws.Range("E20").Value2 = "Ester"
$ws.Range("F20").Value2 = "Guilherme"
$ws.Range("G20").Formula = "=DATE(YEAR(TODAY()),9,2)"
$ws.Range("H20").Formula = "=DATE(YEAR(TODAY()),9,6)"

# --- Row 21 (new): Script Banco de Dados, now Atrasado ---
$ws.Range("B21").Value2 = "Script Banco de Dados"
$ws.Range("C21").Value2 = "Atrasado"
$ws.Range("D21").Value2 = "Semana 2-B"
$ws.Range("E21").Value2 = "Ester"
$ws.Range("G21").Formula = "=DATE(YEAR(TODAY()),9,2)"
$ws.Range("H21").Formula = "=DATE(YEAR(TODAY()),9,6)"

# --- Column width adjustments (approximate, engine-quantized) ---
$ws.Range("B1").ColumnWidth = 65.54
$ws.Range("F1").ColumnWidth = 20.54
